$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.044.18'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.787.35'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.46%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.81'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('E6').Value = '  +1.20%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.46'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.14'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.281'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('E11').Value = '  -2.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0929'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.045.55'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.51'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +11.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.783.03'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.83%  '
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.052.41'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('E18').Value = '  -3.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.52'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '253.05'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0743'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.48'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.27'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.12%  '
$ws.Range('E25').Value = '  -2.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.01'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.58'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.03'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.65%  '
$ws.Range('E29').Value = '  -2.18%  '
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.81'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0517'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.84'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.452.78'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.24%  '
$ws.Range('E37').Value = '  -0.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.629'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '83.47'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.55%  '
$ws.Range('E41').Value = '  -1.83%  '
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.901'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.35%  '
$ws.Range('E44').Value = '  -3.42%  '
$ws.Range('E45').Value = '  -1.73%  '
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.944.02'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.04%  '
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.93'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.22'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.80%  '
